# Apply "added inpatient total cost variables" edit to the Lookup workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Lookup Table")
$ws2 = $wb.Worksheets.Item("Type and Label")

# --- Update existing row 92 --------------------------------------------
# A92: "itotcost" -> "i_adultandpeds_cost"
$ws1.Range("A92").Value = "i_adultandpeds_cost"
# C92: "0600" -> "00300"
$ws1.Range("C92").Value = "00300"
# D92: "20000" -> "03000"
$ws1.Range("D92").Value = "03000"

# --- New rows 93-104: inpatient total cost variables --------------------
# Each row: A=label, B="C000001", C="00300", D=<line number>, F=10, G=1

$rows = @(
    @{ Row = 93;  A = "i_icu_cost";             D = "03100" },
    @{ Row = 94;  A = "i_ccu_cost";              D = "03200" },
    @{ Row = 95;  A = "i_bicu_cost";             D = "03300" },
    @{ Row = 96;  A = "i_sicu_cost";             D = "03400" },
    @{ Row = 97;  A = "i_otherspecial_cost";     D = "03500" },
    @{ Row = 98;  A = "i_subprovideripf_cost";   D = "04000" },
    @{ Row = 99;  A = "i_subproviderirf_cost";   D = "04100" },
    @{ Row = 100; A = "i_subprovider_cost";      D = "04200" },
    @{ Row = 101; A = "i_nursery_cost";          D = "04301" },
    @{ Row = 102; A = "i_skillednursing_cost";   D = "04400" },
    @{ Row = 103; A = "i_nursing_cost";          D = "04500" },
    @{ Row = 104; A = "i_otherlongterm_cost";    D = "04600" }
)

# First pass: columns A, B, C, F, G for every new row (this introduces the
# new "i_..." shared strings in order, before any of the line-number strings
# used by column D get introduced).
foreach ($r in $rows) {
    $row = $r.Row

    $ws1.Cells.Item($row, 1).Value = $r.A

    $ws1.Cells.Item($row, 2).Value = "C000001"
    $ws1.Cells.Item($row, 2).Style = "Normal"

    $ws1.Cells.Item($row, 3).Value = "00300"

    $ws1.Cells.Item($row, 6).Value = 10
    $ws1.Cells.Item($row, 7).Value = 1
}

# Second pass: column D (line numbers) for every new row. The original
# authoring session filled these slightly out of row order (row 102 before
# row 101), which affects where new shared strings land; replicate that
# order here so the shared string table matches exactly.
$dOrder = @(93, 94, 95, 96, 97, 98, 99, 100, 102, 101, 103, 104)
foreach ($row in $dOrder) {
    $r = $rows | Where-Object { $_.Row -eq $row }

    if ($row -eq 93) {
        $ws1.Cells.Item($row, 4).Value = $r.D
    } else {
        $ws1.Cells.Item($row, 4).NumberFormat = "@"
        $ws1.Cells.Item($row, 4).Font.Color = 0
        $ws1.Cells.Item($row, 4).Value = $r.D
    }
}

# --- Sheet view / selection state ---------------------------------------
# "Type and Label" keeps its selection but is no longer the active tab.
$ws2.Range("C51").Select()

# "Lookup Table" becomes the active sheet/tab with updated zoom & scroll.
$ws1.Activate()
$excel.ActiveWindow.Zoom = 119
$ws1.Range("E102").Select()

$wb.Windows.Item(1).Width = 13360
